# Commit: "set no suggests values"
# For each category's "no suggests" catch-all row on sheet "Hoja1",
# shift the previous upper bound (column E) down into the lower bound
# (column D), and set the upper bound to a very large sentinel value
# (999999999) so the range effectively has no practical upper limit /
# "no suggestion" cap.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$noSuggestRows = @(13, 16, 26, 31, 32)
$newLowerBounds = @{
    13 = 4500001
    16 = 10000001
    26 = 8000001
    31 = 3000000
    32 = 10000001
}
$sentinel = 999999999

foreach ($r in $noSuggestRows) {
    $ws.Cells.Item($r, 4).Value = $newLowerBounds[$r]  # column D
    $ws.Cells.Item($r, 5).Value = $sentinel            # column E
}
